# Insert a new record row at row 188 of the weekly Piña price log, shifting
# the existing rows 188-224 down to 189-225, and populate the new row with
# the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new entire row before row 188; remaining rows shift down.
$ws.Rows.Item(188).Insert()

# Populate the newly inserted row 188 with the new weekly record.
$ws.Cells.Item(188, 1).Value = 5
$ws.Cells.Item(188, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(188, 3).Value = "Maule"
$ws.Cells.Item(188, 4).Value = 44641
$ws.Cells.Item(188, 5).Value = 7
$ws.Cells.Item(188, 6).Value = "Fruta"
$ws.Cells.Item(188, 7).Value = 100108
$ws.Cells.Item(188, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(188, 9).Value = 100108005
$ws.Cells.Item(188, 10).Value = "Piña"
$ws.Cells.Item(188, 11).Value = "Caramelo"
$ws.Cells.Item(188, 12).Value = "Segunda"
$ws.Cells.Item(188, 13).Value = 850
$ws.Cells.Item(188, 14).Value = 16000
$ws.Cells.Item(188, 15).Value = 16000
$ws.Cells.Item(188, 16).Value = 16000
$ws.Cells.Item(188, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(188, 18).Value = "Ecuador"
$ws.Cells.Item(188, 19).Value = 1143
$ws.Cells.Item(188, 20).Value = 14

# Make sure date column D keeps the date number format used throughout the
# sheet (same as the other rows in this block).
$ws.Cells.Item(188, 4).NumberFormat = $ws.Cells.Item(189, 4).NumberFormat
